$d = $word.ActiveDocument

# 1. Fix the date: "08," -> "13," (August 08, 2014 -> August 13, 2014)
$d.Content.Find.Execute("08,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "13,", 2)

# 2. Fix the spelling of "excersises" -> "exercises"
$d.Content.Find.Execute("excersises", $true, $false, $false, $false, $false,
                         $true, 1, $false, "exercises", 2)

# 3. Fix the spelling of "seperated" -> "separated"
$d.Content.Find.Execute("seperated", $true, $false, $false, $false, $false,
                         $true, 1, $false, "separated", 2)
